$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 111, shifting existing rows 111-123 down to 112-124.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new weekly data point.
$ws.Cells.Item(111, 1).Value = 7
$ws.Cells.Item(111, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(111, 3).Value = "Ñuble"
$ws.Cells.Item(111, 4).Value = 44449
$ws.Cells.Item(111, 5).Value = 16
$ws.Cells.Item(111, 6).Value = 100112032
$ws.Cells.Item(111, 7).Value = "Zapallo italiano"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 100
$ws.Cells.Item(111, 11).Value = 16000
$ws.Cells.Item(111, 12).Value = 17000
$ws.Cells.Item(111, 13).Value = 16500
$ws.Cells.Item(111, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(111, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(111, 16).Value = 330
$ws.Cells.Item(111, 17).Value = 50
$ws.Cells.Item(111, 18).Value = "Hortaliza"

# Match the date-number style used by the other rows in column D.
$ws.Cells.Item(111, 4).NumberFormat = $ws.Cells.Item(112, 4).NumberFormat
